# "Fruta / hortaliza, semanal"
#
# A new weekly sample for Brócoli at "Macroferia Regional de Talca" is
# inserted at the top of the date-ordered block (row 331); the previously
# existing rows 331-350 shift down one row to 332-351, so row 351 is a
# brand-new row that ends up holding the values that used to live in row
# 350. Row 330 (and everything above it) is untouched.
#
# Columns A,B,C,E,F,G,H,N,Q,R are constant for this whole sub-block, so
# only D (Fecha), I (Calidad), J (Volumen), K/L/M (precios), O (Origen)
# and P (Precio $/Kg) actually vary row to row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$constA = 5
$constB = "Macroferia Regional de Talca"
$constC = "Maule"
$constE = 7
$constF = 100112023
$constG = "Brócoli"
$constH = "Sin especificar"
$constN = '$/unidad'
$constQ = 1
$constR = "Hortaliza"

# Columns: RowNum, D(Fecha), I(Calidad), J(Volumen), K(Precio minimo),
#          L(Precio maximo), M(Precio promedio ponderado), O(Origen),
#          P(Precio $/Kg)
$rows = @(
    @(331, 44746, "Primera", 5000, 800, 800, 800, "Región del Maule", 800),
    @(332, 44690, "Primera", 5000, 700, 700, 700, "Región del Maule", 700),
    @(333, 44631, "Primera", 3000, 500, 500, 500, "Región del Maule", 500),
    @(334, 44235, "Primera", 3000, 600, 600, 600, "Región del Maule", 600),
    @(335, 44307, "Primera", 5000, 400, 400, 400, "Región del Maule", 400),
    @(336, 44672, "Primera", 5000, 600, 600, 600, "Región del Maule", 600),
    @(337, 44344, "Primera", 4000, 600, 600, 600, "Región del Maule", 600),
    @(338, 44707, "Primera", 3000, 1000, 1000, 1000, "Región del Maule", 1000),
    @(339, 44265, "Primera", 3000, 700, 700, 700, "Región del Maule", 700),
    @(340, 44421, "Segunda", 3000, 500, 500, 500, "Región del Maule", 500),
    @(341, 44215, "Primera", 2000, 500, 500, 500, "Región del Maule", 500),
    @(342, 44215, "Segunda", 2000, 300, 300, 300, "Región del Maule", 300),
    @(343, 44566, "Primera", 4000, 500, 500, 500, "Región del Maule", 500),
    @(344, 44637, "Primera", 5000, 400, 400, 400, "Región del Maule", 400),
    @(345, 44483, "Primera", 4000, 800, 800, 800, "Región Metropolitana", 800),
    @(346, 44663, "Primera", 5000, 700, 700, 700, "Región del Maule", 700),
    @(347, 44188, "Primera", 3000, 500, 500, 500, "Región del Maule", 500),
    @(348, 44187, "Primera", 3000, 450, 450, 450, "Región del Maule", 450),
    @(349, 44519, "Primera", 5000, 500, 500, 500, "Región del Maule", 500),
    @(350, 44231, "Primera", 3000, 600, 600, 600, "Región del Maule", 600),
    @(351, 44194, "Primera", 3000, 500, 500, 500, "Región del Maule", 500)
)

foreach ($r in $rows) {
    $rowNum = $r[0]

    $ws.Cells.Item($rowNum, 1).Value2 = $constA
    $ws.Cells.Item($rowNum, 2).Value2 = $constB
    $ws.Cells.Item($rowNum, 3).Value2 = $constC

    $ws.Cells.Item($rowNum, 4).Value2 = $r[1]
    $ws.Cells.Item($rowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($rowNum, 5).Value2 = $constE
    $ws.Cells.Item($rowNum, 6).Value2 = $constF
    $ws.Cells.Item($rowNum, 7).Value2 = $constG
    $ws.Cells.Item($rowNum, 8).Value2 = $constH

    $ws.Cells.Item($rowNum, 9).Value2 = $r[2]
    $ws.Cells.Item($rowNum, 10).Value2 = $r[3]
    $ws.Cells.Item($rowNum, 11).Value2 = $r[4]
    $ws.Cells.Item($rowNum, 12).Value2 = $r[5]
    $ws.Cells.Item($rowNum, 13).Value2 = $r[6]

    $ws.Cells.Item($rowNum, 14).Value2 = $constN
    $ws.Cells.Item($rowNum, 15).Value2 = $r[7]
    $ws.Cells.Item($rowNum, 16).Value2 = $r[8]
    $ws.Cells.Item($rowNum, 17).Value2 = $constQ
    $ws.Cells.Item($rowNum, 18).Value2 = $constR
}

"Updated rows 331-351"
